$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend "A" to each name in column A (rows 1-8)
for ($r = 1; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    $cell.Value = "A" + $current
}

# Update the selection to D8 (single cell), as last used/active cell
$ws.Range("D8").Select()
